# Apply cryptos price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.047.85"
$ws.Range("D3").Value = "'2.953.85"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'594.77"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'148.52"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'2.952.84"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "  +4.04%  "
$ws.Range("E11").Value = "  +6.66%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("D14").Value = "'32.82"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'3.444.82"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'63.017.25"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "'6.70"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'2.952.68"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'442.66"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "'81.06"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").Value = "'11.76"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'7.27"
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +16.32%  "
$ws.Range("D33").Value = "'26.43"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "'3.14"
$ws.Range("E37").Value = "  +5.48%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'49.72"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'8.51"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "'38.87"
$ws.Range("E44").Value = "  -7.73%  "
$ws.Range("D45").Value = "'135.46"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "'2.693.20"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "'360.71"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "'22.85"
$ws.Range("E51").Value = "  -3.10%  "
